$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1922
$ws.Range("J51").Value = 1970.3334
$ws.Range("L51").Value = 1970.3334
$ws.Range("N51").Value = -2938.3334

$ws.Range("H99").Value = 412.77777
$ws.Range("I99").Value = 412.77777
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1238.33331
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H137").Value = 1410.2041
$ws.Range("I137").Value = 1182.6285
$ws.Range("J137").Value = 1979.1428
$ws.Range("K137").Value = 3547.8855
$ws.Range("L137").Value = 5937.428400000001
$ws.Range("M137").Value = -997.8855000000003
$ws.Range("N137").Value = -11037.4284

$ws.Range("H141").Value = 3783.1345
$ws.Range("I141").Value = 1765.5
$ws.Range("J141").Value = 19251.666
$ws.Range("K141").Value = 5296.5
$ws.Range("L141").Value = 57754.99800000001
$ws.Range("M141").Value = -116.5
$ws.Range("N141").Value = -68114.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 20000
$ws.Range("J8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("N8").Value = -20288

$ws.Range("H13").Value = 49866.668

$ws.Range("H32").Value = 9457.187
$ws.Range("I32").Value = 10110.193
$ws.Range("J32").Value = 6342.846
$ws.Range("K32").Value = 10110.193
$ws.Range("L32").Value = 6342.846
$ws.Range("M32").Value = -9823.192999999999
$ws.Range("N32").Value = -6916.846

$ws.Range("H122").Value = 863
$ws.Range("I122").Value = 720.875
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2162.625
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 287.375
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2490.8333
$ws.Range("I134").Value = 2075.6667
$ws.Range("J134").Value = 3072.0667
$ws.Range("K134").Value = 6227.000100000001
$ws.Range("L134").Value = 9216.2001
$ws.Range("M134").Value = -3692.000100000001
$ws.Range("N134").Value = -14286.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 950
$ws.Range("I6").Value = 900
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -787
$ws.Range("N6").Value = -1226

$ws.Range("H12").Value = 10799.6
$ws.Range("I12").Value = 1049.5
$ws.Range("J12").Value = 49800
$ws.Range("K12").Value = 1049.5
$ws.Range("L12").Value = 49800
$ws.Range("M12").Value = -879.5
$ws.Range("N12").Value = -50140

$ws.Range("H82").Value = 29750
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 29750
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -30472

$ws.Range("H85").Value = 29750
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 29750
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32246

$ws.Range("H107").Value = 282.9375
$ws.Range("I107").Value = 281.76923
$ws.Range("J107").Value = 288
$ws.Range("K107").Value = 281.76923
$ws.Range("L107").Value = 288
$ws.Range("M107").Value = 1638.23077
$ws.Range("N107").Value = -4128

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").ClearContents()

$ws.Range("H132").Value = 423963.9
$ws.Range("I132").Value = 521143.5
$ws.Range("J132").Value = 2852.3333
$ws.Range("K132").Value = 1563430.5
$ws.Range("L132").Value = 8556.999899999999
$ws.Range("M132").Value = -1560900.5
$ws.Range("N132").Value = -13616.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 751.8946999999999
$ws.Range("I122").Value = 705.125
$ws.Range("J122").Value = 1001.3333
$ws.Range("K122").Value = 6346.125
$ws.Range("L122").Value = 9011.9997
$ws.Range("M122").Value = -3896.125
$ws.Range("N122").Value = -13911.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 36181.453
$ws.Range("J51").Value = 36181.453
$ws.Range("L51").Value = 36181.453
$ws.Range("N51").Value = -37199.453

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2943.1614
$ws.Range("I7").Value = 2887.8235
$ws.Range("J7").Value = 3010.3572
$ws.Range("K7").Value = 2887.8235
$ws.Range("L7").Value = 3010.3572
$ws.Range("M7").Value = -2775.8235
$ws.Range("N7").Value = -3234.3572

$ws.Range("H68").Value = 1892.1666
$ws.Range("I68").Value = 1643.7142
$ws.Range("J68").Value = 2240
$ws.Range("K68").Value = 1643.7142
$ws.Range("L68").Value = 2240
$ws.Range("M68").Value = -894.7141999999999
$ws.Range("N68").Value = -3738

$ws.Range("H71").Value = 1892.1666
$ws.Range("I71").Value = 1643.7142
$ws.Range("J71").Value = 2240
$ws.Range("K71").Value = 8218.571
$ws.Range("L71").Value = 11200
$ws.Range("M71").Value = -4474.571
$ws.Range("N71").Value = -18688

$ws.Range("H74").Value = 20822.2
$ws.Range("I74").Value = 14555.5
$ws.Range("K74").Value = 14555.5
$ws.Range("M74").Value = -13557.5

$ws.Range("H77").Value = 20822.2
$ws.Range("I77").Value = 14555.5
$ws.Range("K77").Value = 43666.5
$ws.Range("M77").Value = -38674.5

$ws.Range("H126").Value = 2943.1614
$ws.Range("I126").Value = 2887.8235
$ws.Range("J126").Value = 3010.3572
$ws.Range("K126").Value = 8663.470499999999
$ws.Range("L126").Value = 9031.071599999999
$ws.Range("M126").Value = -6193.470499999999
$ws.Range("N126").Value = -13971.0716

$ws.Range("H132").Value = 3010.3333
$ws.Range("I132").Value = 2970.72
$ws.Range("J132").Value = 3100.3635
$ws.Range("K132").Value = 8912.16
$ws.Range("L132").Value = 9301.0905
$ws.Range("M132").Value = -6382.16
$ws.Range("N132").Value = -14361.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4668.6665
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4668.6665
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5916.6665

$ws.Range("H65").Value = 4668.6665
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4668.6665
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -29583.3325

$ws.Range("H82").Value = 31350.25
$ws.Range("J82").Value = 31350.25
$ws.Range("L82").Value = 31350.25
$ws.Range("N82").Value = -32116.25

$ws.Range("H85").Value = 31350.25
$ws.Range("J85").Value = 31350.25
$ws.Range("L85").Value = 31350.25
$ws.Range("N85").Value = -34002.25

$ws.Range("H132").Value = 1307.6296
$ws.Range("I132").Value = 1110.3684
$ws.Range("J132").Value = 1776.125
$ws.Range("K132").Value = 3331.1052
$ws.Range("L132").Value = 5328.375
$ws.Range("M132").Value = -801.1052
$ws.Range("N132").Value = -10388.375

$ws.Range("H135").Value = 43803.75
$ws.Range("J135").Value = 43803.75
$ws.Range("L135").Value = 43803.75
$ws.Range("N135").Value = -53943.75

$ws.Range("H136").Value = 1960
$ws.Range("I136").Value = 2200.125
$ws.Range("J136").Value = 999.5
$ws.Range("K136").Value = 6600.375
$ws.Range("L136").Value = 2998.5
$ws.Range("M136").Value = -4050.375
$ws.Range("N136").Value = -8098.5
